$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A ("Grupo") has its group label only on the first row of each
# group, leaving the rest blank. Fill those blanks down with the group
# label above them (like Excel's Go To Special > Blanks + Fill Down).
$lastRow = $ws.UsedRange.Rows.Count

$lastValue = $null
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($null -ne $val -and $val -ne "") {
        $lastValue = $val
    } elseif ($null -ne $lastValue) {
        $cell.Value2 = $lastValue
    }
}
